$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "LobbyChatMessage"
$ws.Range("C6").Value = "Contiene un messaggio testuale per della chat della lobby, e le informazioni del mittente"

$ws.Range("A6").Select()
